# Apply the "finish parsing Utility" edit:
#  1. Update re-parsed compressor/heat-exchanger coefficients and equipment-cost
#     totals on the CAPCOST sheet (Q8 no longer has a value).
#  2. Add a new UTILITY sheet summarising utility consumption per unit.

$wb = $excel.ActiveWorkbook
$capcost = $wb.Worksheets.Item("CAPCOST")

# --- CAPCOST: re-parsed K1/K2/K3 coefficients for the compressor/HTX columns ---
$capcost.Range("I2").Value = 2.2891
$capcost.Range("J2").Value = 2.2891
$capcost.Range("L2").Value = 2.2891
$capcost.Range("N2").Value = 2.2891
$capcost.Range("S2").Value = 2.2891
$capcost.Range("I3").Value = 1.3604
$capcost.Range("J3").Value = 1.3604
$capcost.Range("L3").Value = 1.3604
$capcost.Range("N3").Value = 1.3604
$capcost.Range("S3").Value = 1.3604
$capcost.Range("I4").Value = -0.1027
$capcost.Range("J4").Value = -0.1027
$capcost.Range("L4").Value = -0.1027
$capcost.Range("N4").Value = -0.1027
$capcost.Range("S4").Value = -0.1027

# --- CAPCOST: re-parsed EQUIPMENT COST row (row 8) ---
$capcost.Range("C8").Value = 1837738.887525055
$capcost.Range("H8").Value = 1976709.767700679
$capcost.Range("I8").Value = 1481215.826683687
$capcost.Range("J8").Value = 994189.8037161557
$capcost.Range("L8").Value = 1221426.706326984
$capcost.Range("M8").Value = 4232444.44902235
$capcost.Range("N8").Value = 1416467.388385196
$capcost.Range("O8").Value = 1470768.43577163
$capcost.Range("S8").Value = 640555.8597020187

# Q8 no longer has a parsed equipment cost -> blank it out
$capcost.Range("Q8").Value = ""

# --- Add the new UTILITY sheet (after CAPCOST) ---
$utility = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$utility.Name = "UTILITY"

# Seed the header row + the two label cells (A2/A3) with the bold/bordered
# "header" style used throughout the workbook (style of CAPCOST!B1 / CAPCOST!A2),
# then overwrite the values/content.
$headerCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")
foreach ($col in $headerCols) {
    $capcost.Range("B1").Copy($utility.Range($col + "1"))
}
$capcost.Range("A2").Copy($utility.Range("A2"))
$capcost.Range("A2").Copy($utility.Range("A3"))

# Row 1: stream/equipment-name headers
$utility.Range("B1").Value = "H-COMP-1"
$utility.Range("C1").Value = "H-COMP-2"
$utility.Range("D1").Value = "H-HTX-1"
$utility.Range("E1").Value = "HTX-02"
$utility.Range("F1").Value = "HTX-03"
$utility.Range("G1").Value = "HTX-04"
$utility.Range("H1").Value = "HTX-05"
$utility.Range("I1").Value = "N-COMP-1"
$utility.Range("J1").Value = "N-COMP-2"
$utility.Range("K1").Value = "N-COMP-3"
$utility.Range("L1").Value = "N-HTX-1"
$utility.Range("M1").Value = "N-HTX-2"
$utility.Range("N1").Value = "NH3FLASH"
$utility.Range("O1").Value = "REACT-1"
$utility.Range("P1").Value = "REACT-2"
$utility.Range("Q1").Value = "REACT-3"

# Row 2: utility-type labels (A2 = type index 0)
$utility.Range("A2").Value = 0
$utility.Range("B2").Value = "ELECTRICITY UTILITY"
$utility.Range("C2").Value = "ELECTRICITY UTILITY"
$utility.Range("D2").Value = "COOLING UTILITY"
$utility.Range("E2").Value = "HOT UTILITY"
$utility.Range("F2").Value = "ELECTRICITY UTILITY"
$utility.Range("G2").Value = "COOLING UTILITY"
$utility.Range("H2").Value = "ELECTRICITY UTILITY"
$utility.Range("I2").Value = "ELECTRICITY UTILITY"
$utility.Range("J2").Value = "ELECTRICITY UTILITY"
$utility.Range("K2").Value = "ELECTRICITY UTILITY"
$utility.Range("L2").Value = "COOLING UTILITY"
$utility.Range("M2").Value = "COOLING UTILITY"
$utility.Range("N2").Value = "HOT UTILITY"
$utility.Range("O2").Value = "HOT UTILITY"
$utility.Range("P2").Value = "HOT UTILITY"
$utility.Range("Q2").Value = "HOT UTILITY"

# Row 3: utility consumption values (A3 = type index 1)
$utility.Range("A3").Value = 1
$utility.Range("B3").Value = 4153.8584
$utility.Range("C3").Value = 3865.109
$utility.Range("D3").Value = 609760
$utility.Range("E3").Value = 1348.87976
$utility.Range("F3").Value = 1070.8652
$utility.Range("G3").Value = 2168900
$utility.Range("H3").Value = 2284.7394
$utility.Range("I3").Value = 2238.2037
$utility.Range("J3").Value = 3059.7322
$utility.Range("K3").Value = 1194.7998
$utility.Range("L3").Value = 372530
$utility.Range("M3").Value = 540000
$utility.Range("N3").Value = -70.37069600000001
$utility.Range("O3").Value = -8882.631999999998
$utility.Range("P3").Value = -5950.903200000001
$utility.Range("Q3").Value = -3961.578560000001

Write-Output "UTILITY sheet created; CAPCOST re-parsed values applied."
